$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume(1h) data scraped this run.
# D-column "Price" values are plain numeric-looking text (e.g. "36.277.86",
# "0.370"); writing them straight through Value2 lets Excel infer a Number
# type and silently reformat/round them (dropping trailing zeros, re-parsing
# the dotted thousands groups, etc). Force the cell to Text, write the
# literal string, then restore the default "Normal" style so no stray
# number-format is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "36.277.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.013.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -0.80%  "

$ws.Range("E4").Value2 = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "252.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +3.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -2.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "62.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +12.66%  "

$ws.Range("E8").Value2 = "  -0.07%  "

$ws.Range("E9").Value2 = "  -5.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.370"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +2.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0744"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +0.92%  "

$ws.Range("E12").Value2 = "  -1.53%  "

$ws.Range("E13").Value2 = "  +1.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "14.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +5.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "2.307.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -1.13%  "

$ws.Range("E16").Value2 = "  +1.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "19.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +15.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "1.998.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -2.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "36.201.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "72.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.0₃0860"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "234.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -0.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +19.93%  "

$ws.Range("E25").Value2 = "  +0.05%  "

$ws.Range("E26").Value2 = "  -1.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +4.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "163.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "19.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -0.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.117"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +35.69%  "

$ws.Range("E31").Value2 = "  +0.20%  "

$ws.Range("E32").Value2 = "  +4.25%  "

$ws.Range("E33").Value2 = "  -0.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "4.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +4.23%  "

$ws.Range("B35").Value2 = "Hedera"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.0606"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +2.03%  "

$ws.Range("B36").Value2 = "LidoDAOToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +14.21%  "

$ws.Range("E37").Value2 = "  -0.14%  "

$ws.Range("E38").Value2 = "  -0.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "5.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +18.07%  "

$ws.Range("E40").Value2 = "  +13.90%  "

$ws.Range("E41").Value2 = "  +1.65%  "

$ws.Range("E42").Value2 = "  +1.76%  "

$ws.Range("E43").Value2 = "  +1.42%  "

$ws.Range("E44").Value2 = "  +3.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "16.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +7.47%  "

$ws.Range("E46").Value2 = "  +6.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "94.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +2.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.428.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +5.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +16.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "47.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +4.71%  "
